# Update meeting minutes for 11 03 25
# Fills in previously-blank "What you have done so far" / "What you want to
# do next" cells for several meeting rows, and corrects one existing entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (meeting 45702, Bhavjot) - wording correction
$ws.Range("D15").Value = "- Completed the introduction for the initial draft report"

# Row 20 (meeting 45705, Mo) - was blank
$ws.Range("D20").Value = "N/A"
$ws.Range("E20").Value = "N/A"

# Row 30 (meeting 45708, Mo) - was blank
$ws.Range("D30").Value = "- Look over colleted datasets once all have been compiled."
$ws.Range("E30").Value = "- Provide feedback on datasets and recommend which datasets use for further analysis"

# Row 35 (meeting 45713, Bhavjot) - was blank
$ws.Range("D35").Value = "-Reviewed Feedback given from Lecturers"
$ws.Range("E35").Value = "- Planned Next steps (When to do presentation & Cleaning datasets for EDA)"

# Row 45 (meeting 45727, Bhavjot) - was blank
$ws.Range("D45").Value = "-EDA / Cleaning a dataset as well as provide data visualisations for analysis`n- Discussed clustering model`n- Began Analysis of EDA (In relation to the draft project report for the client)"
$ws.Range("E45").Value = "- Carry on with EDA and cleaning`n- Get started on the presentation`n- Begin Analysis of EDA"

# Row 46 (Mo) - was blank
$ws.Range("D46").Value = "- Discussed clustering model"
$ws.Range("E46").Value = "- start looking at creating the  clustering model"

# Row 47 (Mourad) - was blank
$ws.Range("D47").Value = "-EDA / Cleaning`n- Discussed clustering model"
$ws.Range("E47").Value = "-Finish resole cleaning and assit Rawad if required"

# Row 48 (Bhavjot) - was blank
$ws.Range("D48").Value = "- Logged feedback for meeting with Phil`n- Discussed clustering model"
$ws.Range("E48").Value = "- Change up the report to include revised work tasks"

# Row 49 (Kel) - was blank
$ws.Range("D49").Value = "- setup Git repository , uploaded all documents so far , and shared with group`n-Submitted initial draft report"
$ws.Range("E49").Value = "- push meeting minutes to GIT `n- Review Mourads code `n-Assist others where needed"
